$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "64.654.59"
$ws.Range("E2").Value = "  -2.78%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.178.36"
$ws.Range("E3").Value = "  -4.55%  "

# Row 4 - TetherUSD
Set-TextValue "D4" "1.00"
$ws.Range("E4").Value = "  -0.04%  "

# Row 5 - BNB
Set-TextValue "D5" "573.27"
$ws.Range("E5").Value = "  -2.26%  "

# Row 6 - Solana
Set-TextValue "D6" "169.87"
$ws.Range("E6").Value = "  -7.04%  "

# Row 7 - XRP
Set-TextValue "D7" "0.607"
$ws.Range("E7").Value = "  -6.38%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.22%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "3.193.60"
$ws.Range("E9").Value = "  -4.03%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -5.15%  "

# Row 11 - Toncoin
Set-TextValue "D11" "6.83"
$ws.Range("E11").Value = "  +0.18%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  -3.02%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "3.730.73"
$ws.Range("E13").Value = "  -4.80%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  -1.69%  "

# Row 15 - WrappedBTC
$ws.Range("D15").Value = "64.687.25"
$ws.Range("E15").Value = "  -2.75%  "

# Row 16 - Avalanche
$ws.Range("E16").Value = "  -4.14%  "

# Row 17 - ShibaInu
Set-TextValue "D17" "0.0000158"
$ws.Range("E17").Value = "  -3.69%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "3.185.36"
$ws.Range("E18").Value = "  -3.58%  "

# Row 19 - BitcoinCash
Set-TextValue "D19" "418.16"
$ws.Range("E19").Value = "  -1.78%  "

# Row 20 - Chainlink
Set-TextValue "D20" "12.98"
$ws.Range("E20").Value = "  -1.39%  "

# Row 21 - Polkadot
$ws.Range("E21").Value = "  -3.56%  "

# Row 22 - Uniswap
Set-TextValue "D22" "7.16"
$ws.Range("E22").Value = "  -3.24%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  -0.09%  "

# Row 24 - LEO
Set-TextValue "D24" "5.69"
$ws.Range("E24").Value = "  +0.07%  "

# Row 25 - Litecoin
Set-TextValue "D25" "70.16"
$ws.Range("E25").Value = "  -2.49%  "

# Row 26 - Kaspa
Set-TextValue "D26" "0.204"
$ws.Range("E26").Value = "  -0.87%  "

# Row 27 - Polygon
$ws.Range("E27").Value = "  -3.07%  "

# Row 28 - PEPE
$ws.Range("E28").Value = "  -9.10%  "

# Row 29 - InternetComputer(DFINITY)
$ws.Range("E29").Value = "  -1.92%  "

# Row 30 - Binance-PegBSC-USD
Set-TextValue "D30" "0.998"
$ws.Range("E30").Value = "  -0.03%  "

# Row 31 - PancakeSwap
Set-TextValue "D31" "1.83"
$ws.Range("E31").Value = "  -4.99%  "

# Row 32 - EthereumClassic
Set-TextValue "D32" "21.78"
$ws.Range("E32").Value = "  -2.99%  "

# Row 33 - USDe
$ws.Range("E33").Value = "  -0.08%  "

# Row 34 - NEARProtocol
Set-TextValue "D34" "5.10"
$ws.Range("E34").Value = "  -2.47%  "

# Row 35 - Aptos
Set-TextValue "D35" "6.38"
$ws.Range("E35").Value = "  -3.94%  "

# Row 36 - Fetch.AI
$ws.Range("E36").Value = "  -4.58%  "

# Row 37 - Monero
Set-TextValue "D37" "155.97"
$ws.Range("E37").Value = "  -2.78%  "

# Row 38 - ImmutableX
$ws.Range("E38").Value = "  -5.29%  "

# Row 39 - Stacks
$ws.Range("E39").Value = "  -5.45%  "

# Row 40 - Maker
$ws.Range("D40").Value = "2.704.18"
$ws.Range("E40").Value = "  -5.52%  "

# Row 41 - Filecoin
Set-TextValue "D41" "4.24"
$ws.Range("E41").Value = "  -2.36%  "

# Row 42 - EnergySwap
$ws.Range("E42").Value = "  -8.29%  "

# Row 43 - OKB
Set-TextValue "D43" "39.13"
$ws.Range("E43").Value = "  -1.63%  "

# Row 44 - Mantle
Set-TextValue "D44" "0.719"
$ws.Range("E44").Value = "  -5.66%  "

# Row 45 - Hedera
Set-TextValue "D45" "0.0623"
$ws.Range("E45").Value = "  -6.13%  "

# Row 46 - RenderToken
Set-TextValue "D46" "5.63"
$ws.Range("E46").Value = "  -6.25%  "

# Row 47 - VeChain
$ws.Range("E47").Value = "  -3.08%  "

# Row 48 - InjectiveProtocol
Set-TextValue "D48" "21.65"
$ws.Range("E48").Value = "  -6.92%  "

# Row 49 - Bittensor
Set-TextValue "D49" "291.32"
$ws.Range("E49").Value = "  -7.21%  "

# Row 50 - Stellar
Set-TextValue "D50" "0.0994"

# Row 51 - dogwifhat -> FirstDigitalUSD
$ws.Range("B51").Value = "FirstDigitalUSD"
$ws.Range("C51").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D51" "0.998"
$ws.Range("E51").Value = "  -0.23%  "
